$d = $word.ActiveDocument

function Add-QAParagraph {
    param(
        [string]$Text,
        [int]$Level
    )
    $endRange = $d.Content
    $endRange.Collapse(0)
    $endRange.InsertParagraphAfter()
    $endRange.Collapse(0)
    $newPara = $d.Paragraphs.Last
    $newPara.Range.Text = $Text
    $newPara.Range.ListFormat.ListLevelNumber = $Level + 1
}

Add-QAParagraph "Which gives the percentage of accurate result" 0
Add-QAParagraph "Recall percentage" 1
Add-QAParagraph "Which score helps to check the quality of model" 0
Add-QAParagraph "Precision" 1
Add-QAParagraph "How can we find the overall performance of the model" 0
Add-QAParagraph "By checking F1 score" 1
Add-QAParagraph "How to find a good model?" 0
Add-QAParagraph "Recall percentage and F1 score" 1
Add-QAParagraph "If recall and precision scores are confusing which one should we check" 0
Add-QAParagraph "F1 score" 1

Write-Output "Added $($d.Paragraphs.Count) total paragraphs"
